# Insert a new data row for "Feria Lagunitas de Puerto Montt - Acelga" at row 94,
# pushing the existing rows 94..174 down to 95..175 (the sheet's dimension grows
# from A1:R174 to A1:R175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(94).Insert()

$ws.Range("A94").Value = 4
$ws.Range("B94").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C94").Value = "Los Lagos"
$ws.Range("D94").Value = 44705
$ws.Range("E94").Value = 10
$ws.Range("F94").Value = 100112009
$ws.Range("G94").Value = "Acelga"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 90
$ws.Range("K94").Value = 12000
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = 12000
$ws.Range("N94").Value = "$/docena de atados (12 kilos)"
$ws.Range("O94").Value = "Región de La Araucanía"
$ws.Range("P94").Value = 1000
$ws.Range("Q94").Value = 12
$ws.Range("R94").Value = "Hortaliza"
